$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C5) from 2023-09-14 to 2023-09-15
# (serial date 45183 -> 45184), one day later, for each of the 4 data rows.
$ws.Range("C2:C5").Value = 45184
